$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Save" header in H1, reusing the same formatting as the other
# header cells (copy G1's formatting onto H1 so it shares the existing
# style entry rather than minting a new one).
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Add the numeric "Save" value in H2
$ws.Range("H2").Value = 1
